# Corrección errores muestreo de datos
#
# Applies the changes described by the commit:
#  - Makes "Modelo de Dominio Anemico" (sheet 1) the active/selected tab
#    instead of "Objetos de dominio" (sheet 2).
#  - On "Notificación" sheet: fixes the F-column concatenation formulas to
#    build a "Tipo-Destinatario-Mensaje" string (instead of "numero Tipo"),
#    wraps text + grows row 3, and nudges a few column widths.
#  - On "Cliente" sheet: nudges a few column widths.
#  - Moves the remembered selection on a couple of sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet references (by tab position, matching workbook.xml <sheets> order)
# ---------------------------------------------------------------------
$wsModelo       = $wb.Worksheets.Item(1)   # Modelo de Dominio Anemico
$wsObjetos      = $wb.Worksheets.Item(2)   # Objetos de dominio
$wsNotificacion = $wb.Worksheets.Item(3)   # Notificación
$wsCliente      = $wb.Worksheets.Item(4)   # Cliente

# ---------------------------------------------------------------------
# Notificación sheet: correct the F column formulas
# ---------------------------------------------------------------------
$wsNotificacion.Range("F2").Formula = "=+B2&""-""&C2&""-""&D2"
$wsNotificacion.Range("F3").Formula = "=+B3&""-""&C3&""-""&D3"
$wsNotificacion.Range("F4").Formula = "=+B4&""-""&C4&""-""&D4"

# F3 gets wrap text (new style: yellow fill + border + wrapText) and its
# row grows to fit the longer wrapped message.
$wsNotificacion.Range("F3").WrapText = $true
$wsNotificacion.Rows.Item(3).RowHeight = 45

# Column width touch-ups on Notificación
$wsNotificacion.Columns.Item(2).ColumnWidth = 15.0
$wsNotificacion.Columns.Item(5).ColumnWidth = 11.666666666666666
$wsNotificacion.Columns.Item(6).ColumnWidth = 72.16666666666667

# Restore the remembered selection on Notificación
$wsNotificacion.Range("D13").Select() | Out-Null

# ---------------------------------------------------------------------
# Cliente sheet: column width touch-ups + remembered selection
# ---------------------------------------------------------------------
$wsCliente.Columns.Item(4).ColumnWidth = 9.5
$wsCliente.Columns.Item(5).ColumnWidth = 26.666666666666668
$wsCliente.Columns.Item(6).ColumnWidth = 17.666666666666668
$wsCliente.Columns.Item(7).ColumnWidth = 27.0

$wsCliente.Range("C16").Select() | Out-Null

# ---------------------------------------------------------------------
# Make "Modelo de Dominio Anemico" the active tab (was "Objetos de dominio")
# ---------------------------------------------------------------------
$wsModelo.Activate()
